# Applies the "Add files via upload" data refresh to the Optical_Power sheet:
#  - Updates a few OT (Orden de Trabajo) codes that had been "Pendiente ADM"
#  - Removes the resolved claim that used to be on row 38 (AMENABAR 3517 / Caso 7432)
#  - Appends two brand-new claims reported on 10/15/2025 (Caso 4146 and Caso 3578)
#
# NOTE: several columns in this sheet store numeric-looking data (case numbers,
# dates typed as plain text, OT codes with leading zeros) as literal TEXT, not
# as Excel numbers/dates. Assigning a numeric-looking string straight to
# .Value would make Excel auto-convert it (stripping leading zeros / turning
# dates into serials), so for those cells we briefly force a Text number
# format, write the value, then clear the format again so the cell keeps its
# plain/default styling while remaining text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) Row 19: the OT code has come back from the admin, fill it in.
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("E19") "00995307"

# ---------------------------------------------------------------------------
# 2) The AMENABAR 3517 claim (old row 38 / Caso 7432) has been closed out and
#    removed entirely; deleting the row shifts every following row up by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(38).Delete()

# ---------------------------------------------------------------------------
# 3) Fill in the OT codes that came back from admin for the rows that shifted
#    up into positions 38-43 (row 39 already carried a real OT code over and
#    needs no change).
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("E38") "00995424"
Set-TextValue $ws.Range("E40") "00995591"
Set-TextValue $ws.Range("E41") "00995740"
Set-TextValue $ws.Range("E42") "00995838"
Set-TextValue $ws.Range("E43") "00996054"

# ---------------------------------------------------------------------------
# 4) Append the two newly reported claims as rows 44 and 45.
# ---------------------------------------------------------------------------
Set-TextValue $ws.Range("A44") "4146"
Set-TextValue $ws.Range("B44") "10/15/2025"
Set-TextValue $ws.Range("C44") "VIRGILIO 1332"
$ws.Range("D44").Value = 10
Set-TextValue $ws.Range("E44") "Pendiente ADM"
Set-TextValue $ws.Range("F44") "Optical Power"
Set-TextValue $ws.Range("G44") "Pendiente"
Set-TextValue $ws.Range("H44") "Cable cortado y en panza"
$ws.Range("I44").Value = 1
Set-TextValue $ws.Range("J44") '{"direccionesNormalizadas": [{"altura": 1332, "cod_calle": 23073, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.510480", "y": "-34.626427"}, "direccion": "VIRGILIO 1332, CABA", "nombre_calle": "VIRGILIO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K44").Value = -58.510480
$ws.Range("L44").Value = -34.626427
Set-TextValue $ws.Range("M44") "Devoto"
Set-TextValue $ws.Range("N44") "Capital Norte"

Set-TextValue $ws.Range("A45") "3578"
Set-TextValue $ws.Range("B45") "10/15/2025"
Set-TextValue $ws.Range("C45") "MIRANDA 3786"
$ws.Range("D45").Value = 11
Set-TextValue $ws.Range("E45") "Pendiente ADM"
Set-TextValue $ws.Range("F45") "Optical Power"
Set-TextValue $ws.Range("G45") "Pendiente"
Set-TextValue $ws.Range("H45") "Caja de empalme colgando"
$ws.Range("I45").Value = 1
Set-TextValue $ws.Range("J45") '{"direccionesNormalizadas": [{"altura": 3786, "cod_calle": 13088, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.493850", "y": "-34.614267"}, "direccion": "MIRANDA 3786, CABA", "nombre_calle": "MIRANDA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K45").Value = -58.493850
$ws.Range("L45").Value = -34.614267
Set-TextValue $ws.Range("M45") "Devoto"
Set-TextValue $ws.Range("N45") "Capital Norte"
